$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the data table (rows 204-207) by copying the formatting of the last
# existing data row (203) down across the new range, then overwrite the
# values with the new observations (Sep-2018 .. Dec-2018).
$ws.Range("A203:D203").Copy($ws.Range("A204:D207"))

$newRows = @(
    @(43344, -0.7, 5.0999999999999996, -9.4),
    @(43374, -1.3, 3.9, -9.1999999999999993),
    @(43405, -19.600000000000001, -23.7, -13.5),
    @(43435, -8.3000000000000007, -11.1, -4.0999999999999996)
)

$startRow = 204
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Move the selection/view to the newly added last cell, like the saved file.
$ws.Range("A188").Select()
$excel.ActiveWindow.ScrollRow = 188
$ws.Range("D207").Select()
